{"js": "// Replace the 25 two-digit-divided-by-one-digit expressions in the\n// worksheet table with their new values, cell by cell (row, col, old, new).\n// Using (row, col) coordinates (rather than a blind document-wide find/\n// replace) avoids ambiguity because some \"before\" values repeat\n// (e.g. \"93\u00f74=\" appears twice but maps to two different results).\n\nconst replacements = [\n  { row: 0, col: 0, oldText: \"19\u00f72=\", newText: \"92\u00f76=\" },\n  { row: 0, col: 1, oldText: \"72\u00f79=\", newText: \"83\u00f76=\" },\n  { row: 0, col: 2, oldText: \"80\u00f79=\", newText: \"40\u00f74=\" },\n  { row: 0, col: 3, oldText: \"11\u00f75=\", newText: \"76\u00f76=\" },\n  { row: 0, col: 4, oldText: \"51\u00f79=\", newText: \"58\u00f77=\" },\n\n  { row: 4, col: 0, oldText: \"62\u00f73=\", newText: \"22\u00f74=\" },\n  { row: 4, col: 1, oldText: \"93\u00f74=\", newText: \"77\u00f76=\" },\n  { row: 4, col: 2, oldText: \"37\u00f73=\", newText: \"85\u00f74=\" },\n  { row: 4, col: 3, oldText: \"86\u00f74=\", newText: \"99\u00f72=\" },\n  { row: 4, col: 4, oldText: \"56\u00f74=\", newText: \"90\u00f76=\" },\n\n  { row: 8, col: 0, oldText: \"64\u00f74=\", newText: \"48\u00f74=\" },\n  { row: 8, col: 1, oldText: \"93\u00f74=\", newText: \"97\u00f72=\" },\n  { row: 8, col: 2, oldText: \"61\u00f76=\", newText: \"74\u00f72=\" },\n  { row: 8, col: 3, oldText: \"72\u00f73=\", newText: \"75\u00f76=\" },\n  { row: 8, col: 4, oldText: \"59\u00f79=\", newText: \"20\u00f73=\" },\n\n  { row: 12, col: 0, oldText: \"46\u00f79=\", newText: \"97\u00f76=\" },\n  { row: 12, col: 1, oldText: \"30\u00f78=\", newText: \"22\u00f73=\" },\n  { row: 12, col: 2, oldText: \"16\u00f75=\", newText: \"35\u00f75=\" },\n  { row: 12, col: 3, oldText: \"14\u00f72=\", newText: \"98\u00f75=\" },\n  { row: 12, col: 4, oldText: \"34\u00f78=\", newText: \"99\u00f75=\" },\n\n  { row: 16, col: 0, oldText: \"56\u00f78=\", newText: \"82\u00f73=\" },\n  { row: 16, col: 1, oldText: \"25\u00f76=\", newText: \"23\u00f79=\" },\n  { row: 16, col: 2, oldText: \"73\u00f73=\", newText: \"39\u00f72=\" },\n  { row: 16, col: 3, oldText: \"99\u00f79=\", newText: \"30\u00f77=\" },\n  { row: 16, col: 4, oldText: \"45\u00f75=\", newText: \"91\u00f76=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const { row, col, oldText, newText } of replacements) {\n  const cell = table.getCell(row, col);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // Replace in place so the existing run formatting (font, size,\n    // paragraph alignment, etc.) is preserved.\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: if for some reason the expected text wasn't found\n    // (e.g. already updated), just overwrite the whole cell body.\n    cell.body.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit-divided-by-one-digit expressions in the\n# worksheet table with their new values, cell by cell (row, col, old, new).\n# Using (row, col) coordinates (rather than a blind document-wide find/\n# replace) avoids ambiguity because some \"before\" values repeat\n# (e.g. \"93\u00f74=\" appears twice but maps to two different results).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{Row=1;  Col=1; Old=\"19\u00f72=\"; New=\"92\u00f76=\"},\n    @{Row=1;  Col=2; Old=\"72\u00f79=\"; New=\"83\u00f76=\"},\n    @{Row=1;  Col=3; Old=\"80\u00f79=\"; New=\"40\u00f74=\"},\n    @{Row=1;  Col=4; Old=\"11\u00f75=\"; New=\"76\u00f76=\"},\n    @{Row=1;  Col=5; Old=\"51\u00f79=\"; New=\"58\u00f77=\"},\n\n    @{Row=5;  Col=1; Old=\"62\u00f73=\"; New=\"22\u00f74=\"},\n    @{Row=5;  Col=2; Old=\"93\u00f74=\"; New=\"77\u00f76=\"},\n    @{Row=5;  Col=3; Old=\"37\u00f73=\"; New=\"85\u00f74=\"},\n    @{Row=5;  Col=4; Old=\"86\u00f74=\"; New=\"99\u00f72=\"},\n    @{Row=5;  Col=5; Old=\"56\u00f74=\"; New=\"90\u00f76=\"},\n\n    @{Row=9;  Col=1; Old=\"64\u00f74=\"; New=\"48\u00f74=\"},\n    @{Row=9;  Col=2; Old=\"93\u00f74=\"; New=\"97\u00f72=\"},\n    @{Row=9;  Col=3; Old=\"61\u00f76=\"; New=\"74\u00f72=\"},\n    @{Row=9;  Col=4; Old=\"72\u00f73=\"; New=\"75\u00f76=\"},\n    @{Row=9;  Col=5; Old=\"59\u00f79=\"; New=\"20\u00f73=\"},\n\n    @{Row=13; Col=1; Old=\"46\u00f79=\"; New=\"97\u00f76=\"},\n    @{Row=13; Col=2; Old=\"30\u00f78=\"; New=\"22\u00f73=\"},\n    @{Row=13; Col=3; Old=\"16\u00f75=\"; New=\"35\u00f75=\"},\n    @{Row=13; Col=4; Old=\"14\u00f72=\"; New=\"98\u00f75=\"},\n    @{Row=13; Col=5; Old=\"34\u00f78=\"; New=\"99\u00f75=\"},\n\n    @{Row=17; Col=1; Old=\"56\u00f78=\"; New=\"82\u00f73=\"},\n    @{Row=17; Col=2; Old=\"25\u00f76=\"; New=\"23\u00f79=\"},\n    @{Row=17; Col=3; Old=\"73\u00f73=\"; New=\"39\u00f72=\"},\n    @{Row=17; Col=4; Old=\"99\u00f79=\"; New=\"30\u00f77=\"},\n    @{Row=17; Col=5; Old=\"45\u00f75=\"; New=\"91\u00f76=\"}\n)\n\n# NOTE: wdReplaceAll (2), even when scoped to a single cell's Range, has\n# been observed to corrupt sibling-cell lookups in this host (a later\n# Cell() call can return stale/wrong text). wdReplaceOne (1) is both\n# sufficient (each cell holds exactly one match) and safe.\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    $rng = $cell.Range\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($item.Old, $false, $false, $false, $false, $false, $true, 1, $false, $item.New, 1) | Out-Null\n}\n"}
